$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 'Jake Windham_20251202_124710'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("C16").Value = 'Jake Windham'
$ws.Range("D16").Value = 20
$ws.Range("E16").Value = 'Male'
$ws.Range("F16").Value = '2025-12-02 12:47:10'
$ws.Range("G16").Value = '{
  "portion": 0.8,
  "diet": 1.0,
  "salt": 0.8,
  "fat": 0.8,
  "natural": 1.0,
  "convenience": 0.4,
  "price": 1.0
}'
$ws.Range("H16").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = '0.571'
$ws.Range("J16").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K16").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = '0.483'
$ws.Range("M16").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N16").Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range("O16").NumberFormat = "@"
$ws.Range("O16").Value = '0.459'
$ws.Range("P16").Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range("Q16").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = '0.662'
$ws.Range("S16").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("T16").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("U16").NumberFormat = "@"
$ws.Range("U16").Value = '0.605'
$ws.Range("V16").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("W16").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("X16").NumberFormat = "@"
$ws.Range("X16").Value = '0.602'
$ws.Range("Y16").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z16").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = '0.690'
$ws.Range("AB16").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC16").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AD16").NumberFormat = "@"
$ws.Range("AD16").Value = '0.557'
$ws.Range("AE16").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AF16").Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AG16").NumberFormat = "@"
$ws.Range("AG16").Value = '0.535'
$ws.Range("AH16").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

$ws.Rows.Item(16).AutoFit()